$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.2895481763322857
$ws.Range("E2").Value = 0.2895481763322857

# Row 3
$ws.Range("D3").Value = 0.0001042542288175873
$ws.Range("E3").Value = 0.0001042542288175873

# Row 4
$ws.Range("C4").Value = $false
$ws.Range("D4").Value = 0.9998974032622709
$ws.Range("E4").Value = 0.9998974032622709

# Row 5
$ws.Range("D5").Value = 0.9997289538032235
$ws.Range("E5").Value = 0.9997289538032235

# Row 6
$ws.Range("D6").Value = 0.9999999999999984
$ws.Range("E6").Value = 0.9999999999999984

# Row 7
$ws.Range("F7").Value = 8.643099784851074
$ws.Range("G7").Value = 0.5
